# OCMTimelineReport.xlsx - "Commit Reports Fixed 03/09/2021"
#
# 1. Queries!G2  - fix QueryDrillGridOne SQL: drop "AST.AgentID as [Agent ID],"
#    from the select list and change the final sort from desc to Asc.
# 2. AdvanceSearch!I2 / J2 - update the sample search strings used by the
#    advanced-search demo row (Search String1 / Search String2).
# 3. Re-create the navigation/selection state that was left behind when the
#    author saved the file (active sheet + selected cell on several tabs).

$wb = $excel.ActiveWorkbook

# --- 1. Fix the QueryDrillGridOne SQL text -------------------------------
$wsQueries = $wb.Worksheets.Item("Queries")

$newQueryDrillGridOne = @"
SELECT distinct ISNULL(A.FirstName,'')+' '+ ISNULL(A.LastName,'') AS [Agent Name],
AST.StationID as [Station ID], 
AST.[Status],AST.StatusType as [Status Type], 
format(dbo.VARCHARTODATETIME(AST.[TimeStamp]),'dd/MM/yyyy HH:mm:ss') as [TimeStamp]   
FROM  [dbo].[AGT_Agent_StatusTrack] AST  with(nolock)
left join AGT_Agent a on a.AvayaLoginID=AST.agentid
WHERE format(dbo.VARCHARTODATETIME( AST.[TimeStamp]),'yyyy-MM-dd HH:mm:ss') >= 'ReportBeforeDate'
AND format(dbo.VARCHARTODATETIME( AST.[TimeStamp]),'yyyy-MM-dd HH:mm:ss') <='ReportAfterDate' 
AND AgentID like 'AgentIdCapturedFromUI' ORDER BY [TimeStamp] Asc

"@
$wsQueries.Range("G2").Value = $newQueryDrillGridOne

# --- 2. Update the AdvanceSearch demo row search strings -----------------
$wsAdvance = $wb.Worksheets.Item("AdvanceSearch")
$wsAdvance.Range("I2").Value = "'6666"
$wsAdvance.Range("J2").Value = "'Ab"

# --- 3. Restore the sheet/cell selection state ---------------------------
$wsQueries.Activate()
$wsQueries.Range("A2").Select()

$wsAdvance.Activate()
$wsAdvance.Range("J3").Select()

$wsShowInNewPageDateRange = $wb.Worksheets.Item("ShowInNewPageDateRange")
$wsShowInNewPageDateRange.Activate()
$wsShowInNewPageDateRange.Range("B28").Select()

$wsShow = $wb.Worksheets.Item("Show")
$wsShow.Activate()
$wsShow.Range("C12").Select()
